$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 374, pushing the existing rows 374-402 down to 375-403
$ws.Rows(374).Insert()

# Populate the newly inserted row with the new weekly record
$ws.Range("A374").Value = 10
$ws.Range("B374").Value = "Vega Modelo de Temuco"
$ws.Range("C374").Value = "La Araucanía"
$ws.Range("D374").Value = 45106
$ws.Range("E374").Value = 9
$ws.Range("F374").Value = 100112052
$ws.Range("G374").Value = "Albahaca"
$ws.Range("H374").Value = "Sin especificar"
$ws.Range("I374").Value = "Primera"
$ws.Range("J374").Value = 110
$ws.Range("K374").Value = 6000
$ws.Range("L374").Value = 6000
$ws.Range("M374").Value = 6000
$ws.Range("N374").Value = "$/paquete"
$ws.Range("O374").Value = "Región de Arica y Parinacota"
$ws.Range("P374").Value = 6000
$ws.Range("Q374").Value = 1
$ws.Range("R374").Value = "Hortaliza"
